$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 8 (shifts old row 8 "Upload" down to row 9, etc.)
$ws.Rows.Item(8).Insert()

# Fill the new row 8 with the "Force" entry, matching formatting/style of neighboring rows
$ws.Cells.Item(8, 1).Value = "Force"
for ($col = 2; $col -le 10; $col++) {
    $ws.Cells.Item(8, $col).Value = $false
}

# Copy style (s="3") from row 7/9 onto new row 8 cells
$ws.Range("A7:J7").Copy()
$ws.Range("A8:J8").PasteSpecial(-4122)
$ws.Cells.Item(8, 1).Value = "Force"
for ($col = 2; $col -le 10; $col++) {
    $ws.Cells.Item(8, $col).Value = $false
}

$ws.Rows.Item(8).RowHeight = 13.5
$ws.Rows.Item(7).RowHeight = 13.5
$ws.Rows.Item(9).RowHeight = 13.5
